# Fixed Weapon Spec: Maul bug, updated skill changes notes.
# Adds two new class-skill-change entries (Wizard, 2E Psionicist) to the
# bottom of the "Athasian Skill Changes" table on Sheet1, following the
# same 2-row-per-class layout (New Class Skills / Cross-Class Skills)
# used by every other class already in the sheet, separated by a blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Wizard
$ws.Range("A51").Value = "Wizard"
$ws.Range("B51").Value = "Bluff"
$ws.Range("C51").Value = "Swim"
$ws.Range("B52").Value = "Disguise?"
$ws.Range("C52").Value = "Survival"

# (row 53 intentionally left blank, matching the separator pattern)

# 2E Psionicist
$ws.Range("A54").Value = "2E Psionicist"
$ws.Range("C54").Value = "Swim"
$ws.Range("C55").Value = "Survival"

# Move the active selection to just past the new last row, matching the
# author's saved view state.
$ws.Range("C56").Select()
